# Repro: add explicit "currency_movements" tab, fix typing hints and other minor issues
$wb = $excel.ActiveWorkbook

# --- leave the previously active sheet (currency_conversions) with a stray
#     selection at E33, matching where the author's cursor ended up, and
#     make sure it is no longer the tab that is marked as selected.
$cc = $wb.Worksheets.Item("currency_conversions")
$cc.Activate() | Out-Null
$cc.Range("E33").Select() | Out-Null

# --- add the new sheet as the last tab, right after currency_conversions
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "currency_movements"

$dateFmt = "yyyy\-mm\-dd;@"

# --- header row (bold, like the other sheets in the workbook)
$ws.Range("A1:E1").Font.Bold = $true
$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "buy_date"
$ws.Range("C1").Value = "amount"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "comment"

# --- row 2
$ws.Range("A2:B2").NumberFormat = $dateFmt
$ws.Range("A2").Value = (Get-Date -Year 2024 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B2").Value = (Get-Date -Year 2022 -Month 10 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = "USD"

# --- row 3
$ws.Range("A3:B3").NumberFormat = $dateFmt
$ws.Range("A3").Value = (Get-Date -Year 2024 -Month 1 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B3").Value = (Get-Date -Year 2024 -Month 1 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = "EUR"

# --- row 4 (date columns get the darker/explicit black font color)
$ws.Range("A4:B4").NumberFormat = $dateFmt
$ws.Range("A4:B4").Font.Color = 0
$ws.Range("A4").Value = (Get-Date -Year 2024 -Month 1 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B4").Value = (Get-Date -Year 2024 -Month 1 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C4").Value = -100
$ws.Range("D4").Value = "EUR"

# --- row 5 (date columns get the darker/explicit black font color)
$ws.Range("A5:B5").NumberFormat = $dateFmt
$ws.Range("A5:B5").Font.Color = 0
$ws.Range("A5").Value = (Get-Date -Year 2024 -Month 1 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B5").Value = (Get-Date -Year 2024 -Month 1 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C5").Value = -100
$ws.Range("D5").Value = "USD"

# --- make the new sheet the active / selected tab
$ws.Activate() | Out-Null
$ws.Range("A1").Select() | Out-Null
